# Actualización automática 2025-10-15 17:30:09
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (ventas por grupo de producto)
# ---------------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("L32").Value = 89.56
$wsGrupo.Range("D35").Value = 177.06
$wsGrupo.Range("M35").Value = 2848.12
$wsGrupo.Range("M47").Value = 1356

$wsGrupo.Range("D60").Value = "1 de 58"
$wsGrupo.Range("L60").Value = "2 de 58"
$wsGrupo.Range("M60").Value = "6 de 58"

# ---------------------------------------------------------------------------
# Sheet "VENTA MENSUAL" (ventas por mes)
# ---------------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F32").Value = 89.56
$wsMensual.Range("F35").Value = 3025.18
$wsMensual.Range("F47").Value = 1356
$wsMensual.Range("F60").Value = 20257.12

# Column D width changed from 14 to 13 characters
$wsMensual.Columns.Item(4).ColumnWidth = 12.17

# ---------------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL" (resumen de cumplimiento)
# ---------------------------------------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumplimiento.Range("D3").Value = 177.06
$wsCumplimiento.Range("E3").Value = 5327.558903862629
$wsCumplimiento.Range("F3").Value = 0.03216571448311449

$wsCumplimiento.Range("D11").Value = 880.34
$wsCumplimiento.Range("E11").Value = 2626.32949822329
$wsCumplimiento.Range("F11").Value = 0.2510473258018865

$wsCumplimiento.Range("D12").Value = 15998.3
$wsCumplimiento.Range("E12").Value = 22543.95
$wsCumplimiento.Range("F12").Value = 0.415084744663324

$wsCumplimiento.Range("D14").Value = 20257.12
$wsCumplimiento.Range("E14").Value = 34226.64774946896
$wsCumplimiento.Range("F14").Value = 0.3718010122418055
